# Apply cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.915.51'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.623.25'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'213.67"
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -2.33%  '
$ws.Range("E9").Value = '  -3.67%  '
$ws.Range("E10").Value = '  -6.81%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").Value = '1.849.39'
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").Value = '1.621.74'
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = "'0.524"
$ws.Range("E15").Value = '  -3.74%  '
$ws.Range("D16").Value = '25.908.73'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = "'61.12"
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("D18").Value = '0.0₃0733'
$ws.Range("E18").Value = '  -3.87%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = "'191.50"
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("E21").Value = '  -3.02%  '
$ws.Range("D22").Value = "'9.57"
$ws.Range("E22").Value = '  -3.54%  '
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").Value = "'143.73"
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").Value = "'0.0483"
$ws.Range("E31").Value = '  -2.52%  '
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("E33").Value = '  -5.52%  '
$ws.Range("E34").Value = '  -3.13%  '
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D36").Value = '1.118.33'
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("D39").Value = "'0.517"
$ws.Range("E39").Value = '  -4.26%  '
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D41").Value = "'97.76"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").Value = "'0.767"
$ws.Range("E42").Value = '  -3.72%  '
$ws.Range("D43").Value = '1.759.20'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("E44").Value = '  -5.78%  '
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").Value = "'0.0530"
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").Value = "'54.35"
$ws.Range("E47").Value = '  -3.80%  '
$ws.Range("D48").Value = "'1.46"
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = '  -4.01%  '
